$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'22.18"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'5.355"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.05924"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'3.395"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'6.393"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.8125"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.9613"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.1431"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.07409"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.03467"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").Value = "'4.059"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'0.09403"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'0.001598"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.04830"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.0005912"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("D19").Value = "'0.006139"
$ws.Range("D19").Style = "Normal"
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D20").Value = "'0.004085"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "19HotbitTokenHTB"
$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D21").Value = "'0.0009874"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "20BitKanKAN"
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D22").Value = "'0.00009702"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "21NitroExNTX"
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").Value = "'3.741"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "22LEOLEO"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D24").Value = "'2.165"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "23BTSETokenBTSE"
$ws.Range("B25").Value = "BitpandaEcosystemToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D25").Value = "'0.3268"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "24BitpandaEcosystemTokenBEST"
$ws.Range("B26").Value = "ProBitToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D26").Value = "'0.1333"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "25ProBitTokenPROB"
$ws.Range("B27").Value = "UpBots"
$ws.Range("C27").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D27").Value = "'0.0002462"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "26UpBotsUBXT"
$ws.Range("D40").Value = "'0.03942"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.006492"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
$ws.Range("D42").Value = "'0.1071"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'0.002701"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'0.005901"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005284"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Value = "'0.6602"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.04705"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
